# Se agregan los scripts 0728/0729 (y relacionados 0730-0732) a la clase Tests_MiPortal
# Extends the DEC_07xx data table: inserts two new rows (so the existing
# trailer rows 13-21 shift down to 15-23) and fills the freed rows 13-17
# with five new records cloned from row 12's layout, only varying column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing trailing rows (old 13-21) down by 2, to rows 15-23.
$ws.Rows("13:14").Insert() | Out-Null

# Clone row 12's full formatting + values into the 5 freshly inserted rows.
$ws.Range("A12:J12").Copy() | Out-Null
for ($r = 13; $r -le 17; $r++) {
    $ws.Range("A$r`:J$r").PasteSpecial(-4104) | Out-Null
}

# Column A gets a distinct script id per new row; B, C and D:J keep the
# values copied from row 12 (18092588-0 / sebA$1357 / SIN_DATO).
$newLabels = @("DEC_0728", "DEC_0729", "DEC_0730", "DEC_0731", "DEC_0732")
for ($i = 0; $i -lt $newLabels.Length; $i++) {
    $row = 13 + $i
    $ws.Cells.Item($row, 1).Value = $newLabels[$i]
}

# Match the author's final selection/active cell.
$ws.Range("D15").Select() | Out-Null
